$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.287.58"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.682.88"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'218.33"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'0.5263"
$ws.Range("E6").Value = "  +2.58%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("D9").Value = "'0.06421"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'22.04"
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("D11").Value = "'0.07485"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.549"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.680.21"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "'0.5812"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "'0.000008491"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "'64.31"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "26.336.20"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "'4.923"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").Value = "'1.007"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "'10.87"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "'188.94"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'144.30"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "'7.707"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("D28").Value = "'0.06629"
$ws.Range("E28").Value = "  +12.38%  "
$ws.Range("D29").Value = "'1.346"
$ws.Range("E29").Value = "  +4.93%  "
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").Value = "'3.576"
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("D32").Value = "'3.560"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").Value = "'0.6194"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("D37").Value = "'2.700"
$ws.Range("E37").Value = "  +1.72%  "
$ws.Range("D38").Value = "'6.407"
$ws.Range("E38").Value = "  +5.73%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.110.49"
$ws.Range("E39").Value = "  +2.85%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01623"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").Value = "'0.8758"
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").Value = "'1.015"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").Value = "'100.55"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").Value = "1.831.17"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "'0.00000000114"
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("D46").Value = "'56.79"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'8.157"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").Value = "'6.038"
$ws.Range("E51").Value = "  +2.85%  "
